$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 137, shifting existing rows 137..211 down to 138..212
$ws.Rows.Item(137).Insert(-4121)

# Populate the new row 137 with this week's data point (weekly Fruta/Hortaliza update)
$ws.Range("A137").Value = 8
$ws.Range("B137").Value = "Terminal La Palmera de La Serena"
$ws.Range("C137").Value = "Coquimbo"
$ws.Range("D137").Value = 44981
$ws.Range("E137").Value = 4
$ws.Range("F137").Value = 100112001
$ws.Range("G137").Value = "Berenjena"
$ws.Range("H137").Value = "Sin especificar"
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 480
$ws.Range("K137").Value = 11000
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = 11500
$ws.Range("N137").Value = "`$/caja 50 unidades"
$ws.Range("O137").Value = "Región de Arica y Parinacota"
$ws.Range("P137").Value = 230
$ws.Range("Q137").Value = 50
$ws.Range("R137").Value = "Hortaliza"
